$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Capture the long / accented text blocks (and a couple of short labels)
#    from their ORIGINAL locations before any structural changes happen.
#    Reading through a PS variable (instead of Write-Output) keeps the
#    original UTF-8 text intact.
# ---------------------------------------------------------------------------
$docente         = $ws.Range("B13").Value2   # "5840560 - Marco Antonio Carvalho Pereira"
$shortSyllabusEn = $ws.Range("B15").Value2   # Short syllabus (EN)
$dataAtivacao    = $ws.Range("B8").Value2    # "01/01/2020"
$syllabusEnLong  = $ws.Range("B17").Value2   # Long syllabus (EN)
$osAlunos        = $ws.Range("B19").Value2   # "Os alunos montarao equipes..."
$avaliacaoTrab   = $ws.Range("B20").Value2   # "Avaliacao dos trabalhos..."
$nfTexto         = $ws.Range("B21").Value2   # "NF = (MF + PR)/ 2 ..."

$lblProgramaResumido = $ws.Range("A14").Value2  # "Programa resumido:"
$lblShortSyllabus    = $ws.Range("A15").Value2  # "Short syllabus:"
$lblPrograma         = $ws.Range("A16").Value2  # "Programa:"
$lblSyllabus         = $ws.Range("A17").Value2  # "Syllabus:"
$lblAvaliacao        = $ws.Range("A18").Value2  # "Avaliacao:"
$lblMetodo           = $ws.Range("A19").Value2  # "Metodo:"
$lblNormaRecup       = $ws.Range("A21").Value2  # "Norma de recuperacao:"
$lblBibliografia     = $ws.Range("A22").Value2  # "Bibliografia:"

# ---------------------------------------------------------------------------
# 2. Drop one row so the grid goes from 22 rows down to 21 (everything we
#    still need has already been captured above).
# ---------------------------------------------------------------------------
$ws.Rows.Item(22).Delete()

# ---------------------------------------------------------------------------
# 3. Row 10 (Objetivos:) now carries the "docente" text instead of the long
#    objectives paragraph.
# ---------------------------------------------------------------------------
$ws.Range("B10").Value = $docente
$ws.Range("C10").Value = $docente

# ---------------------------------------------------------------------------
# 4. Row 13 becomes "Programa resumido:" / "Semestral" (60pt row height).
# ---------------------------------------------------------------------------
$ws.Rows.Item(13).RowHeight = 60
$ws.Range("A13").Value = $lblProgramaResumido
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# ---------------------------------------------------------------------------
# 5. Row 14 becomes "Short syllabus:" / short syllabus EN text (60pt).
# ---------------------------------------------------------------------------
$ws.Rows.Item(14).RowHeight = 60
$ws.Range("A14").Value = $lblShortSyllabus
$ws.Range("B14").Value = $shortSyllabusEn
$ws.Range("C14").Value = $shortSyllabusEn

# ---------------------------------------------------------------------------
# 6. Row 15 becomes "Programa:" / "01/01/2020" (120pt row height).
# ---------------------------------------------------------------------------
$ws.Rows.Item(15).RowHeight = 120
$ws.Range("A15").Value = $lblPrograma
$ws.Range("B15").Value = $dataAtivacao
$ws.Range("C15").Value = $dataAtivacao

# ---------------------------------------------------------------------------
# 7. Row 16 becomes "Syllabus:" / long syllabus EN text (120pt).
# ---------------------------------------------------------------------------
$ws.Rows.Item(16).RowHeight = 120
$ws.Range("A16").Value = $lblSyllabus
$ws.Range("B16").Value = $syllabusEnLong
$ws.Range("C16").Value = $syllabusEnLong

# ---------------------------------------------------------------------------
# 8. Row 17 becomes just "Avaliacao:" (default height, no B/C).
# ---------------------------------------------------------------------------
$ws.Rows.Item(17).RowHeight = 15
$ws.Range("A17").Value = $lblAvaliacao
$ws.Range("B17").ClearContents()
$ws.Range("C17").ClearContents()

# ---------------------------------------------------------------------------
# 9. Row 18 becomes "Metodo:" / docente text again (60pt).
# ---------------------------------------------------------------------------
$ws.Rows.Item(18).RowHeight = 60
$ws.Range("A18").Value = $lblMetodo
$ws.Range("B18").Value = $docente
$ws.Range("C18").Value = $docente

# ---------------------------------------------------------------------------
# 10. Row 19 becomes "Criterio:" (new label) / "Os alunos..." text (60pt).
# ---------------------------------------------------------------------------
$ws.Rows.Item(19).RowHeight = 60
$ws.Range("A19").Value = "Crit" + [char]0x00E9 + "rio:"
$ws.Range("B19").Value = $osAlunos
$ws.Range("C19").Value = $osAlunos

# ---------------------------------------------------------------------------
# 11. Row 20 becomes "Norma de recuperacao:" / "Avaliacao dos trabalhos..." (60pt).
# ---------------------------------------------------------------------------
$ws.Rows.Item(20).RowHeight = 60
$ws.Range("A20").Value = $lblNormaRecup
$ws.Range("B20").Value = $avaliacaoTrab
$ws.Range("C20").Value = $avaliacaoTrab

# ---------------------------------------------------------------------------
# 12. Row 21 becomes "Bibliografia:" / "NF = ..." text (120pt).
# ---------------------------------------------------------------------------
$ws.Rows.Item(21).RowHeight = 120
$ws.Range("A21").Value = $lblBibliografia
$ws.Range("B21").Value = $nfTexto
$ws.Range("C21").Value = $nfTexto
